$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1884.5
$ws.Range("J17").Value = 1884.5
$ws.Range("L17").Value = 5653.5
$ws.Range("N17").Value = -5989.5

$ws.Range("H19").Value = 285.4375
$ws.Range("I19").Value = 255.55556
$ws.Range("J19").Value = 323.85715
$ws.Range("K19").Value = 255.55556
$ws.Range("L19").Value = 323.85715
$ws.Range("M19").Value = -80.55556000000001
$ws.Range("N19").Value = -673.85715

$ws.Range("H64").Value = 37039596
$ws.Range("I64").Value = 76925040
$ws.Range("J64").Value = 3111.4285
$ws.Range("K64").Value = 76925040
$ws.Range("L64").Value = 3111.4285
$ws.Range("M64").Value = -76924792
$ws.Range("N64").Value = -3607.4285

$ws.Range("H67").Value = 37039596
$ws.Range("I67").Value = 76925040
$ws.Range("J67").Value = 3111.4285
$ws.Range("K67").Value = 76925040
$ws.Range("L67").Value = 3111.4285
$ws.Range("M67").Value = -76924182
$ws.Range("N67").Value = -4827.4285

$ws.Range("H112").Value = 1530789.6
$ws.Range("J112").Value = 1927418.9
$ws.Range("L112").Value = 5782256.699999999
$ws.Range("N112").Value = -5784472.699999999

$ws.Range("H116").Value = 3277.2273
$ws.Range("I116").Value = 2822.6
$ws.Range("K116").Value = 2822.6
$ws.Range("M116").Value = 619.4000000000001

$ws.Range("H138").Value = 2156.7046
$ws.Range("I138").Value = 1184.6
$ws.Range("J138").Value = 3435.7896
$ws.Range("K138").Value = 3553.8
$ws.Range("L138").Value = 10307.3688
$ws.Range("M138").Value = 1586.2
$ws.Range("N138").Value = -20587.3688

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 209838.33
$ws.Range("I61").Value = 1328.7646
$ws.Range("J61").Value = 716218.7
$ws.Range("K61").Value = 1328.7646
$ws.Range("L61").Value = 716218.7
$ws.Range("M61").Value = -1116.7646
$ws.Range("N61").Value = -716642.7

$ws.Range("H74").Value = 4866.654
$ws.Range("I74").Value = 873.25
$ws.Range("K74").Value = 873.25
$ws.Range("M74").Value = 0.75

$ws.Range("H77").Value = 4866.654
$ws.Range("I77").Value = 873.25
$ws.Range("K77").Value = 4366.25
$ws.Range("M77").Value = 1.75

$ws.Range("H102").Value = 1525.5555
$ws.Range("I102").Value = 1355
$ws.Range("J102").Value = 1866.6666
$ws.Range("K102").Value = 1355
$ws.Range("L102").Value = 1866.6666
$ws.Range("M102").Value = 267
$ws.Range("N102").Value = -5110.6666

$ws.Range("H130").Value = 24363.572
$ws.Range("J130").Value = 24363.572
$ws.Range("L130").Value = 24363.572
$ws.Range("N130").Value = -34403.572

$ws.Range("H136").Value = 209838.33
$ws.Range("I136").Value = 1328.7646
$ws.Range("J136").Value = 716218.7
$ws.Range("K136").Value = 3986.2938
$ws.Range("L136").Value = 2148656.1
$ws.Range("M136").Value = -1436.2938
$ws.Range("N136").Value = -2153756.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1336.8889
$ws.Range("I134").Value = 1075.8182
$ws.Range("J134").Value = 2485.6
$ws.Range("K134").Value = 3227.4546
$ws.Range("L134").Value = 7456.799999999999
$ws.Range("M134").Value = -692.4546
$ws.Range("N134").Value = -12526.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3979.443
$ws.Range("I31").Value = 1547.25
$ws.Range("J31").Value = 5040.7637
$ws.Range("K31").Value = 1547.25
$ws.Range("L31").Value = 5040.7637
$ws.Range("M31").Value = -1252.25
$ws.Range("N31").Value = -5630.7637

$ws.Range("H34").Value = 3979.443
$ws.Range("I34").Value = 1547.25
$ws.Range("J34").Value = 5040.7637
$ws.Range("K34").Value = 1547.25
$ws.Range("L34").Value = 5040.7637
$ws.Range("M34").Value = -1345.25
$ws.Range("N34").Value = -5444.7637

$ws.Range("H99").Value = 1547.92
$ws.Range("I99").Value = 1319.125
$ws.Range("K99").Value = 1319.125
$ws.Range("M99").Value = 178.875

$ws.Range("H121").Value = 27563
$ws.Range("J121").Value = 27563
$ws.Range("L121").Value = 27563
$ws.Range("N121").Value = -30183

$ws.Range("H122").Value = 951.13043
$ws.Range("I122").Value = 889.3333
$ws.Range("K122").Value = 2667.9999
$ws.Range("M122").Value = -217.9998999999998

$ws.Range("H126").Value = 1547.92
$ws.Range("I126").Value = 1319.125
$ws.Range("K126").Value = 3957.375
$ws.Range("M126").Value = -1487.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1144.8235
$ws.Range("I107").Value = 1192.8572
$ws.Range("J107").Value = 1111.2
$ws.Range("K107").Value = 3578.5716
$ws.Range("L107").Value = 3333.6
$ws.Range("M107").Value = -1658.5716
$ws.Range("N107").Value = -7173.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1701.3043
$ws.Range("I68").Value = 1673.5714
$ws.Range("J68").Value = 1744.4445
$ws.Range("K68").Value = 1673.5714
$ws.Range("L68").Value = 1744.4445
$ws.Range("M68").Value = -924.5714
$ws.Range("N68").Value = -3242.4445

$ws.Range("H71").Value = 1701.3043
$ws.Range("I71").Value = 1673.5714
$ws.Range("J71").Value = 1744.4445
$ws.Range("K71").Value = 8367.857
$ws.Range("L71").Value = 8722.2225
$ws.Range("M71").Value = -4623.857
$ws.Range("N71").Value = -16210.2225

$ws.Range("H82").Value = 1258.8
$ws.Range("I82").Value = 1163.6666
$ws.Range("J82").Value = 1401.5
$ws.Range("K82").Value = 1163.6666
$ws.Range("L82").Value = 1401.5
$ws.Range("M82").Value = -802.6666
$ws.Range("N82").Value = -2123.5

$ws.Range("H85").Value = 1258.8
$ws.Range("I85").Value = 1163.6666
$ws.Range("J85").Value = 1401.5
$ws.Range("K85").Value = 1163.6666
$ws.Range("L85").Value = 1401.5
$ws.Range("M85").Value = 84.33339999999998
$ws.Range("N85").Value = -3897.5

$ws.Range("H132").Value = 10104.909
$ws.Range("I132").Value = 4392.316
$ws.Range("J132").Value = 22874.234
$ws.Range("K132").Value = 13176.948
$ws.Range("L132").Value = 68622.702
$ws.Range("M132").Value = -10646.948
$ws.Range("N132").Value = -73682.702

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 887.40625
$ws.Range("I126").Value = 753.56525
$ws.Range("J126").Value = 1229.4445
$ws.Range("K126").Value = 2260.69575
$ws.Range("L126").Value = 3688.3335
$ws.Range("M126").Value = 209.3042500000001
$ws.Range("N126").Value = -8628.333500000001

